$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 139; existing rows 139:166 shift down to 140:167.
$ws.Rows("139").Insert()

# Populate the newly inserted row 139 with the new price-report record.
$ws.Range("A139").Value = 7
$ws.Range("B139").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C139").Value = "Ñuble"
$ws.Range("D139").Value = 44476
$ws.Range("E139").Value = 16
$ws.Range("F139").Value = 100112009
$ws.Range("G139").Value = "Acelga"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 120
$ws.Range("K139").Value = 350
$ws.Range("L139").Value = 400
$ws.Range("M139").Value = 375
$ws.Range("N139").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O139").Value = "Provincia de Diguillín"
$ws.Range("P139").Value = 375
$ws.Range("Q139").Value = 1
$ws.Range("R139").Value = "Hortaliza"
